# Update absenteeism data rows 2-11 with new values per commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2;  A=72525; B="Dra. Luiza Duarte";        C="Jurídico";               D="Consulta médica";     E=6; F=45091; G=12158.65}
    @{Row=3;  A=79511; B="Davi Luiz Santos";          C="P&D";                    D="Viagem de negócios";  E=1; F=45104; G=3167.51}
    @{Row=4;  A=13622; B="Ana Carolina da Rocha";     C="Atendimento ao Cliente"; D="Problemas pessoais";  E=7; F=45091; G=5868.65}
    @{Row=5;  A=77458; B="Amanda Costa";              C="Recursos Humanos";       D="Problemas pessoais";  E=5; F=45084; G=4120.79}
    @{Row=6;  A=48359; B="Leonardo Ramos";            C="Marketing";              D="Consulta médica";     E=7; F=45090; G=8572.92}
    @{Row=7;  A=47994; B="Dra. Emanuelly Correia";    C="TI";                     D="Viagem de negócios";  E=7; F=45081; G=5559.05}
    @{Row=8;  A=29335; B="Agatha Barbosa";            C="P&D";                    D="Problemas pessoais";  E=6; F=45092; G=9761.83}
    @{Row=9;  A=57370; B="Dra. Emilly Duarte";        C="Operações";              D="Viagem de negócios";  E=3; F=45089; G=2908.59}
    @{Row=10; A=53093; B="Maysa Porto";               C="Atendimento ao Cliente"; D="Problemas pessoais";  E=1; F=45094; G=4752.51}
    @{Row=11; A=38707; B="Sra. Ana Julia Ferreira";   C="TI";                     D="Doença";               E=2; F=45081; G=3365.24}
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
    $ws.Cells.Item($r, 7).Value = $item.G
}

$wb.Save()
